$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set C1 to "Warna" with the same style as B1 (header style)
$ws.Range("C1").Value = "Warna"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Adjust column C width (stored width ends up ColumnWidth + 5/6, so back it out)
$ws.Columns.Item(3).ColumnWidth = 17.1666666666667

# Set the selection to E14
$ws.Range("E14").Select()
